$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45902
$ws.Range("B2").Value = 100.86
$ws.Range("C2").Value = 95
$ws.Range("D2").Value = 87.2
$ws.Range("E2").Value = 81.88
$ws.Range("F2").Value = 80.3
$ws.Range("G2").Value = 80.3
$ws.Range("H2").Value = 81.55
$ws.Range("I2").Value = 95.15000000000001
$ws.Range("J2").Value = 87.2
$ws.Range("K2").Value = 39
$ws.Range("L2").Value = 7.96
$ws.Range("M2").Value = 4.01
$ws.Range("N2").Value = 4.31
$ws.Range("O2").Value = 4.01
$ws.Range("P2").Value = 4.01
$ws.Range("Q2").Value = 4.31
$ws.Range("R2").Value = 4.01
$ws.Range("S2").Value = 14
$ws.Range("T2").Value = 32.89
$ws.Range("U2").Value = 84
$ws.Range("V2").Value = 109.52
$ws.Range("W2").Value = 160
$ws.Range("X2").Value = 116.72
$ws.Range("Y2").Value = 104.22
$ws.Range("Z2").Value = 61.77
$ws.Range("AB2").Value = 122.62
$ws.Range("AD2").Value = 134.76
$ws.Range("AF2").Value = 110.47
$ws.Range("AG2").Value = "9h-18h"
